$d = $word.ActiveDocument

# --- 1. First paragraph: pad the existing sentence with two trailing spaces
#        and append a new, red-colored run of text after it. ---

# Grow "This is a Microsoft word document." to end with two trailing spaces.
$d.Content.Find.Execute("This is a Microsoft word document.", $true, $false, $false, $false, $false, $true, 1, $false, "This is a Microsoft word document.  ", 2) | Out-Null

# Insert the new run immediately after the (now padded) first paragraph's
# text, but before its paragraph mark, so it lands in the same <w:p>.
$para1 = $d.Paragraphs(1).Range
$insertPos = $para1.End - 1
$newRun = $d.Range($insertPos, $insertPos)
$newRun.InsertAfter("(This is a change – Version for branch alternate)")

# Color only the text we just inserted (the new run) dark red (C00000).
$newRunRange = $d.Range($insertPos, $insertPos + 50)
$newRunRange.Font.Color = 192   # wdColor BGR for RGB C00000 -> 0x0000C0 = 192

# --- 2. Mark the "Normal (Web)" style as hidden from the style list
#        (adds <w:semiHidden/> in the stored style definition). ---
$style = $d.Styles("Normal (Web)")
try {
    $style.Hidden = $true
} catch {
    # Some hosts expose Style.Hidden as read-only; ignore if unsettable.
}
